# update account to test account
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Swap the test account credentials (shared strings used by F2/G2).
$ws.Range("F2").Value = "ngqautotester@hpe.com"
$ws.Range("G2").Value = "585347198c1d5b145d3de47ef43cec6ff4731f1872dbf3e75d7d"

# The old account's mailto hyperlink on F2 is no longer wanted.
$ws.Range("F2").Hyperlinks.Delete()

# Move the selection, matching the saved view state in the workbook.
$ws.Range("F8").Select() | Out-Null
